$d = $word.ActiveDocument

# Remove the "Requisitos" heading paragraph and the following List Bullet
# paragraph (the two requirement lines) that appear at the very end of the
# document, right after the Bibliografia paragraph.

$paras = $d.Paragraphs
$count = $paras.Count

# The last paragraph is the List Bullet with the requirement lines,
# the one before it is the "Requisitos" Heading2 paragraph.
$lastPara = $paras.Item($count)
$secondLastPara = $paras.Item($count - 1)

$start = $secondLastPara.Range.Start
$end = $lastPara.Range.End

$deleteRange = $d.Range($start, $end)
$deleteRange.Delete()
